$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Footer "today" date field on the slide master + every slide layout
#    rolled from 11/15/14 to 11/16/14 (the deck was re-saved a day later).
# ---------------------------------------------------------------------------
$m = $p.SlideMaster

$dateTargets = @($m)
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $dateTargets += $m.CustomLayouts.Item($li)
}

foreach ($t in $dateTargets) {
    $shapes = $t.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "11/16/14"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Module numbering bump: every "Module N" label becomes "Module N+1"
#    (slides 1-3), and the longer title on slide 4 gets the same treatment.
# ---------------------------------------------------------------------------

# Slide 1: "Module 1" -> "Module 2"
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(21)
$sh1.TextFrame.TextRange.Text = "Module "
[void]$sh1.TextFrame.TextRange.InsertAfter("2")
$sh1.Width = 84.55370078740158

# Slide 2: "Module 2" -> "Module 3"
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(21)
$sh2.TextFrame.TextRange.Text = "Module "
[void]$sh2.TextFrame.TextRange.InsertAfter("3")
$sh2.Width = 84.55370078740158

# Slide 3: "Module 3" -> "Module 4"
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$sh3.TextFrame.TextRange.Text = "Module "
[void]$sh3.TextFrame.TextRange.InsertAfter("4")
$sh3.Width = 84.55370078740158

# Slide 4: "Module 4 " Rerun Cufflinks in alternative 'modes'" -> "Module 5 " ..."
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(1)
$sh4.TextFrame.TextRange.Text = "Module "
[void]$sh4.TextFrame.TextRange.InsertAfter("5")
[void]$sh4.TextFrame.TextRange.InsertAfter(" ")
[void]$sh4.TextFrame.TextRange.InsertAfter([char]0x2013 + " Rerun Cufflinks in alternative " + [char]0x2018 + "modes" + [char]0x2019)
# box is the same pixel width for "4" and "5" so the autosized box keeps its
# original footprint here (only the boxes on slides 1-3 shrink slightly)
$sh4.Height = 26.62504
